# Append the new daily data row (2025/10/08, 水, 10, 14) as row 77,
# mirroring the existing data rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-looking string ("2025/10/08") that must stay plain
# text (matching every other row in the column), so prefix it with an
# apostrophe the way a user typing it in Excel would, to suppress
# automatic date recognition.
$ws.Cells.Item(77, 1).Value = "'2025/10/08"
$ws.Cells.Item(77, 2).Value = "水"
$ws.Cells.Item(77, 3).Value = 10
$ws.Cells.Item(77, 4).Value = 14
